# TC29_Canine_StudyUBC-AllBreeds_StageOfDisease.xlsx
# - FilesTab query (cell B4, "startup" sheet) gains a `DISTINCT` on its
#   final RETURN clause (and a few stray blank/whitespace-only lines that
#   were sitting in the middle of the query text are cleaned up), so the
#   Files tab report no longer produces duplicate rows.
# - The user's last selection on the sheet moved from C3 to B4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("B4")
$query = [string]$cell.Value2

# Drop the three stray blank / whitespace-only lines that sat between
# "MATCH (p:program)..." and "WITH DISTINCT f, parent, c, demo, diag, s",
# between that line and "WHERE ...", and between "WHERE ..." and
# "OPTIONAL MATCH ...".
$query = $query -replace "(?m)^[ \t]*\r?\n(WITH DISTINCT f, parent, c, demo, diag, s)", '$1'
$query = $query -replace "(?m)(WITH DISTINCT f, parent, c, demo, diag, s)\r?\n\r?\n(WHERE)", "`$1`n`$2"
$query = $query -replace "(?m)(WHERE s\.clinical_study_designation[^\r\n]*)\r?\n\r?\n(OPTIONAL MATCH)", "`$1`n`$2"

# The final RETURN clause becomes RETURN DISTINCT.
$query = $query -replace "(?m)^RETURN[ \t]*\r?\n", "RETURN DISTINCT`n"

$cell.Value2 = $query

# Re-assigning a long wrapped cell's value makes the engine recompute an
# autofit row height; the row used a fixed custom height of 120 before
# and after this edit, so put it back.
$ws.Rows.Item(4).RowHeight = 120

# Reflect where the user's selection ended up after this edit.
$ws.Range("B4").Select() | Out-Null
